$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values on row 2 (Ambiente, URL, Contrasenia, NumPoliza change; Usuario & TIPO_ENDOSO stay the same).
# A2 and E2 keep their original "quote-prefixed text" cell style, so re-enter those
# values with a leading apostrophe to preserve that formatting (and avoid E2's
# numeric-looking value being auto-converted to a number).
# Edit order (E2, A2, B2, D2) matches how the shared-string table was rebuilt in
# the target workbook.
$ws.Range("E2").Value = "'04104016708"
$ws.Range("A2").Value = "'ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B2").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws.Range("D2").Value = "gw"

# Add hyperlink on B2 pointing at the URL value, preserving the existing cell style
$u = $ws.Range("B2").Value2
$styleBefore = $ws.Range("B2").Style
$ws.Hyperlinks.Add($ws.Range("B2"), $u)
$ws.Range("B2").Style = $styleBefore

# Move the active selection to E2
$ws.Range("E2").Select()
